$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that used to sit right after
#    "Deletes an employee" (it will be re-added later, right before
#    the final period of the 3rd LINQ-query bullet).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Replace the whole tail of the 3rd bullet point:
#    ". Select the **project's name**, **start date**, **end date** and **manager name**."
#    with:
#    ". Select each employee's **first name**, **last name** and **manager name**
#      and each of their projects' **name**, **start date**, **end date**."
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute(". Select the project's name, start date, end date and manager name.", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the target sentence to replace"
}

# Wipe the whole matched range, then collapse to its start so we can
# rebuild it run-by-run with the correct bold formatting.
$rng.Text = ""
$rng.Collapse(1)

function Insert-Chunk($range, [string]$text, [bool]$bold) {
    $range.Text = $text
    $range.Bold = [int]$bold
    $range.Collapse(0)
}

Insert-Chunk $rng ". " $false
Insert-Chunk $rng "Select each employee's " $false
Insert-Chunk $rng "first name" $true
Insert-Chunk $rng ", " $false
Insert-Chunk $rng "last name" $true
Insert-Chunk $rng " " $true
Insert-Chunk $rng "and " $false
Insert-Chunk $rng "manager " $true
Insert-Chunk $rng "name " $true
Insert-Chunk $rng "and each of their projects'" $false
Insert-Chunk $rng " name" $true
Insert-Chunk $rng ", " $false
Insert-Chunk $rng "start date" $true
Insert-Chunk $rng ", " $false

# Final boundary: type the trailing "." *first* (while still in the
# plain, non-bold ", " context), then rewind and insert the bold
# "end date" right in front of it. Doing it in this order keeps the
# "." run free of any stray explicit (non-)bold formatting - typing
# "end date" then toggling Bold back off for "." would otherwise
# leave a redundant/empty run-properties element behind.
$rng.Text = "."
$rng.Collapse(1)
$rng.Text = "end date"
$rng.Bold = 1
$rng.Collapse(0)

# ------------------------------------------------------------------
# 3) Re-insert the "_GoBack" bookmark right here - between "end date"
#    and the trailing period.
# ------------------------------------------------------------------
[void]$d.Bookmarks.Add("_GoBack", $rng)

Write-Output "Done"
